$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.899.14"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "1.626.65"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.19"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.521"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.39"
$ws.Range("E8").Value = "  +8.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.259"
$ws.Range("E9").Value = "  +3.02%  "
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0915"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").Value = "1.858.71"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("D13").Value = "1.627.99"
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.576"
$ws.Range("E14").Value = "  +7.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.91"
$ws.Range("E15").Value = "  +4.96%  "
$ws.Range("D16").Value = "29.954.36"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.93"
$ws.Range("E17").Value = "  +17.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.72"
$ws.Range("E18").Value = "  +2.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.64"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").Value = "0.0₃0706"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.14"
$ws.Range("E22").Value = "  +3.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.63"
$ws.Range("E23").Value = "  +4.41%  "
$ws.Range("E24").Value = "  +2.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.38"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.71"
$ws.Range("E26").Value = "  +2.45%  "
$ws.Range("E27").Value = "  +2.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.62"
$ws.Range("E28").Value = "  +3.15%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0490"
$ws.Range("E30").Value = "  +3.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.12"
$ws.Range("E31").Value = "  +5.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  +3.92%  "
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("D34").Value = "1.423.19"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("E35").Value = "  +6.83%  "
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("E37").Value = "  +1.89%  "
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("E39").Value = "  +3.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.558"
$ws.Range("E40").Value = "  +3.59%  "
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("E42").Value = "  +3.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0498"
$ws.Range("E43").Value = "  +2.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "54.30"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.43"
$ws.Range("E45").Value = "  +5.52%  "
$ws.Range("E46").Value = "  +8.06%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.41"
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("D49").Value = "1.766.38"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "88.88"
$ws.Range("E50").Value = "  +2.69%  "
$ws.Range("D51").Value = "0.0₆0109"
$ws.Range("E51").Value = "  +8.25%  "
